# Consumption_Actual_Historical.xlsx update
# Commit message: "Updating the models for the portfolio assets"
#
# The historical-data window rolls forward by 22 days:
#   - rows 2..97   : 30.12.2025 -> 21.01.2026  (date serial 46021 -> 46043)
#   - rows 98..193 : 31.12.2025 -> 22.01.2026  (date serial 46022 -> 46044)
# Column B (Actual Consumption (MW)) is refreshed with the newly fetched
# readings for rows 2..135; rows 136..193 have no reading yet for the new
# day and stay at 0, same as they did for the previous day's tail end.
# Column D ("Lookup") is rebuilt from the (new) date text + the existing
# Quarter index in column C, mirroring how the sheet already encodes it
# (e.g. "30.12.20251" == date "30.12.2025" + quarter "1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 193
$dayShift = 22

# Newly fetched "Actual Consumption (MW)" values for rows 2..135.
# Rows 136..193 are not in this list and are left at 0.
$newConsumption = @(6746,6766,6706,6666,6653,6623,6598,6570,6670,6624,6628,6591,6658,6697,6744,6770,6837,6949,6999,7049,7248,7462,7626,7901,8235,8424,8559,8636,8881,9019,9057,9114,9183,9296,9208,9144,9100,9148,9048,8989,8967,8949,8830,8834,8771,8724,8678,8659,8573,8534,8466,8437,8363,8375,8351,8363,8355,8506,8401,8489,8551,8673,8665,8769,8802,8977,9051,9124,9128,9182,9120,9123,9100,9060,8982,8952,8876,8841,8738,8620,8512,8379,8273,8143,7962,7804,7640,7475,7313,7253,7107,7010,6819,6756,6657,6677,6627,6604,6507,6530,6440,6493,6426,6432,6453,6455,6473,6544,6535,6496,6555,6584,6646,6557,6601,6712,6989,7189,7353,7609,7945,8141,8286,8404,8562,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $idx = $r - $firstRow

    # Shift the timestamp forward by 22 days; the time-of-day fraction is
    # preserved automatically by plain serial-date arithmetic.
    $oldDate = $ws.Cells.Item($r, 1).Value2
    $newDate = $oldDate + $dayShift
    $ws.Cells.Item($r, 1).Value2 = $newDate

    # Refresh the consumption reading for this quarter-hour.
    if ($idx -lt $newConsumption.Length) {
        $newB = $newConsumption[$idx]
    } else {
        $newB = 0
    }
    $ws.Cells.Item($r, 2).Value2 = $newB

    # Rebuild the "Lookup" helper column: "dd.MM.yyyy" of the new date,
    # directly concatenated with the (unchanged) quarter index.
    $quarter = $ws.Cells.Item($r, 3).Value2
    $dateOnly = [DateTime]::FromOADate($newDate)
    $dateText = $dateOnly.ToString("dd.MM.yyyy")
    $ws.Cells.Item($r, 4).Value2 = "$dateText$quarter"
}
